$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.428482174873352
$ws.Range("B1").Value = 3.277657985687256
$ws.Range("C1").Value = 2.301731109619141
$ws.Range("D1").Value = 2.04340672492981
$ws.Range("E1").Value = 1.763611912727356
